$wb = $excel.ActiveWorkbook

# --- "Logs" sheet: append new row 34 with the 6th test mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Cells.Item(34, 1).Value = "Hebben we EcoPro-700 nog op voorraad?"
$ws.Cells.Item(34, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item(34, 3).Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$ws.Cells.Item(34, 4).Value = "Inkoop / Bestellingen"
$ws.Cells.Item(34, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$ws.Cells.Item(34, 6).Value = "2025-08-03 18:23:53"
$ws.Cells.Item(34, 7).Value = "Ja"
$ws.Cells.Item(34, 8).Value = "Ja"
$ws.Cells.Item(34, 9).Value = "Nee"
$ws.Cells.Item(34, 10).Value = "Nee"

# Extend the existing conditional formatting ranges (was *2:*33) so they
# also cover the newly added row 34, keeping the same rules/priorities.
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "33")
    $newRange = $ws.Range($col + "2:" + $col + "34")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- "Dashboard" sheet: bump the "Inkoop / Bestellingen" count from 5 to 6 ---
$wsd = $wb.Worksheets.Item("Dashboard")
$wsd.Cells.Item(5, 2).Value = 6
